$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (my-chapter_name / inline) is unchanged.

# Row 3: paragraph-with-a-really-wide-rule-name -> inline-anchors-in-paragraph
$ws.Range("A3").Value = "my-chapter_name"
$ws.Range("B3").Value = "inline-anchors-in-paragraph"
$ws.Range("C3").Value = "Paragraph with inline anchor and something."
$ws.Range("D3").Value = "[""norm:paragraph:inline-anchors-in-paragraph""]"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""

# Row 4: note_with_2_tags -> inline-anchors-in-tagged-paragraph
$ws.Range("A4").Value = "my-chapter_name"
$ws.Range("B4").Value = "inline-anchors-in-tagged-paragraph"
$ws.Range("C4").Value = "inline anchor"
$ws.Range("D4").Value = "[""norm:paragraph:inline-anchors-in-paragraph:inline-anchors""]"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""

# Row 5: desc1 -> paragraph-with-a-really-wide-rule-name
$ws.Range("A5").Value = "my-chapter_name"
$ws.Range("B5").Value = "paragraph-with-a-really-wide-rule-name"
$ws.Range("C5").Value = "Here's a description.`nIt's got 2 lines.`nParagraph without inline anchors"
$ws.Range("D5").Value = "Description, [""norm:paragraph:no-inline-anchors-in-paragraph""]"

# Row 6: desc2 -> table1 (with formula-looking content in C6)
$ws.Range("A6").Value = "my-chapter_name"
$ws.Range("B6").Value = "table1"
$ws.Range("C6").Formula = "===`n cell with anchor`ncell without anchor`n==="
$ws.Range("D6").Value = "[""norm:table:anchors-in-cells:entire-table""]"

# Row 7: rule_with_newlines -> table2 (with formula-looking content in C7)
$ws.Range("A7").Value = "my-chapter_name"
$ws.Range("B7").Value = "table2"
$ws.Range("C7").Formula = "=""Header 1|Header 2`n==`nCell in column 1, row 1|Cell in column 2, row 1`nCell in column 1, row 2|Cell in column 2, row 2`n==="""
$ws.Range("D7").Value = "[""norm:table:no-anchors-in-cells:entire-table""]"

# Row 8: new row - unordered1
$ws.Range("B8:D8").WrapText = $true
$ws.Range("A8").Value = "my-chapter_name"
$ws.Range("B8").Value = "unordered1"
$ws.Range("C8").Value = "Unordered List with anchors:"
$ws.Range("D8").Value = "[""norm:unordered-list:anchors-in-items:entire-list""]"

# Row 9: new row - note_with_2_tags (re-added, with Kind/Instances)
$ws.Range("B9:D9").WrapText = $true
$ws.Range("A9").Value = "my-chapter_name"
$ws.Range("B9").Value = "note_with_2_tags"
$ws.Range("C9").Value = "One line description`nParagraph 1`nParagraph 3"
$ws.Range("D9").Value = "Description, [""norm:admonition:anchors-in-notes:note1"", ""norm:admonition:anchors-in-notes:note3""]"
$ws.Range("E9").Value = "parameter"
$ws.Range("F9").Value = "MY_PARAMETER"

# Row 10: new row - desc1 (re-added)
$ws.Range("B10:D10").WrapText = $true
$ws.Range("A10").Value = "my-chapter_name"
$ws.Range("B10").Value = "desc1"
$ws.Range("C10").Value = "Description Item 1`nDescription Item 3"
$ws.Range("D10").Value = "[""norm:description-list:anchors-in-items:item1"", ""norm:description-list:anchors-in-items:item3""]"

# Row 11: new row - desc2 (re-added)
$ws.Range("B11:D11").WrapText = $true
$ws.Range("A11").Value = "my-chapter_name"
$ws.Range("B11").Value = "desc2"
$ws.Range("C11").Value = "Description Item 1`nDescription Item 3"
$ws.Range("D11").Value = "[""norm:description-list:anchors-in-items:item1"", ""norm:description-list:anchors-in-items:item3""]"

# Row 12: new row - rule_with_newlines (re-added, description text joined onto one line)
$ws.Range("B12:D12").WrapText = $true
$ws.Range("A12").Value = "my-chapter_name"
$ws.Range("B12").Value = "rule_with_newlines"
$ws.Range("C12").Value = "Here&#8217;s the first line. Here&#8217;s the second line."
$ws.Range("D12").Value = "[""norm:tag_with_newlines""]"

# Resize the table to cover the new data range.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F11"))

# Multi-line values written into wrap-text cells auto-expand the row
# height; reset rows back to the sheet's default (no explicit height),
# matching the original/target formatting.
$ws.Range("A3:F12").EntireRow.AutoFit()
